# Home page loaded test case updated
#
# - Update the existing "User Login with Valid Credentials" result row:
#     Status:         FAILED -> PASSED
#     Execution Time: Tue Mar 25 10:52:53 IST 2025 -> Tue Mar 25 16:11:15 IST 2025
# - Append a new test result row for "Verify Home Page Loads Successfully"
#     Status:         PASSED
#     Execution Time: Tue Mar 25 16:11:15 IST 2025

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (existing test case result)
$ws.Range("B2").Value = "PASSED"
$ws.Range("C2").Value = "Tue Mar 25 16:11:15 IST 2025"

# Add row 3 (new test case result)
$ws.Range("A3").Value = "Verify Home Page Loads Successfully"
$ws.Range("B3").Value = "PASSED"
$ws.Range("C3").Value = "Tue Mar 25 16:11:15 IST 2025"
